# Apply daily-scrape update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (C, D, F, G, H changed; A, B, E unchanged) ---
# Excel's ColumnWidth (characters) gets re-quantized to pixels on save, so the
# stored <col width="..."/> ends up ~0.83 higher than the raw input unless we
# compensate; subtracting 0.85 lands safely inside the pixel bucket that
# serializes back out to the exact target integer width.
$ws.Columns.Item(3).ColumnWidth = 51 - 0.85
$ws.Columns.Item(4).ColumnWidth = 21 - 0.85
$ws.Columns.Item(6).ColumnWidth = 16 - 0.85
$ws.Columns.Item(7).ColumnWidth = 15 - 0.85
$ws.Columns.Item(8).ColumnWidth = 22 - 0.85

# OPPORTUNITY ID (col A) looks numeric but must stay plain text, matching the
# source data. Force text entry by flipping the cell to a text number-format
# before assigning the digit string, then ClearFormats() to drop back to the
# workbook's default (unstyled) cell format without Excel re-coercing the
# stored value back into a number.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Row 2 new data ---
Set-TextValue $ws.Range("A2") "1329433"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1329433"
$ws.Range("C2").Value = "Sales management experience at a trading company"
$ws.Range("D2").Value = "日本、東京都東京"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "4 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Tsuchiya Co., Ltd.,"

# --- Row 3 new data ---
Set-TextValue $ws.Range("A3") "1327775"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327775"
$ws.Range("C3").Value = "Accelerate Romania| Programming Intern"
$ws.Range("D3").Value = "Bucharest, Romania"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "69 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "AQUAsoft"

# --- Row 4 new data ---
Set-TextValue $ws.Range("A4") "1327768"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327768"
$ws.Range("C4").Value = "Accelerate Romania| Business Development Intern"
$ws.Range("D4").Value = "Bucharest, Romania"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "38 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "AQUAsoft"

# --- Remove now-stale rows 5 through 10 ---
$ws.Range("A5:H10").Delete()
